# Scheduled-runner update: refresh market-board price snapshots and
# recompute dependent profit columns across the Kraken_Profits sheets
# (one worksheet per crafting job: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
# Row 8
$ws.Range("H8").Value = 500
$ws.Range("J8").Value = 500
$ws.Range("L8").Value = 1500
$ws.Range("N8").Value = -1778
# Row 17
$ws.Range("H17").Value = 1950
$ws.Range("J17").Value = 1950
$ws.Range("L17").Value = 5850
$ws.Range("N17").Value = -6186
# Row 42
$ws.Range("H42").Value = 51.5
$ws.Range("I42").Value = 1
$ws.Range("J42").Value = 68.333336
$ws.Range("K42").Value = 3
$ws.Range("L42").Value = 205.000008
$ws.Range("M42").Value = 227
$ws.Range("N42").Value = -665.000008
# Row 64
$ws.Range("H64").Value = 3580.4
$ws.Range("I64").Value = 3299.6667
$ws.Range("J64").Value = 4001.5
$ws.Range("K64").Value = 3299.6667
$ws.Range("L64").Value = 4001.5
$ws.Range("M64").Value = -3051.6667
$ws.Range("N64").Value = -4497.5
# Row 67
$ws.Range("H67").Value = 3580.4
$ws.Range("I67").Value = 3299.6667
$ws.Range("J67").Value = 4001.5
$ws.Range("K67").Value = 3299.6667
$ws.Range("L67").Value = 4001.5
$ws.Range("M67").Value = -2441.6667
$ws.Range("N67").Value = -5717.5
# Row 106
$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("M106").ClearContents()
# Row 137
$ws.Range("H137").Value = 4000
$ws.Range("I137").Value = 4000
$ws.Range("K137").Value = 12000
$ws.Range("M137").Value = -9450

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3525.95
$ws.Range("J32").Value = 11595
$ws.Range("L32").Value = 11595
$ws.Range("N32").Value = -12169
# Row 74
$ws.Range("H74").Value = 1550.6
$ws.Range("I74").Value = 1594
$ws.Range("J74").Value = 1268.5
$ws.Range("K74").Value = 1594
$ws.Range("L74").Value = 1268.5
$ws.Range("M74").Value = -720
$ws.Range("N74").Value = -3016.5
# Row 77
$ws.Range("H77").Value = 1550.6
$ws.Range("I77").Value = 1594
$ws.Range("J77").Value = 1268.5
$ws.Range("K77").Value = 7970
$ws.Range("L77").Value = 6342.5
$ws.Range("M77").Value = -3602
$ws.Range("N77").Value = -15078.5
# Row 122
$ws.Range("H122").Value = 10326.096
$ws.Range("I122").Value = 9825.388999999999
$ws.Range("K122").Value = 29476.167
$ws.Range("M122").Value = -27026.167
# Row 135
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("N135").Value = 0
$ws.Range("L135").ClearContents()

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
# Row 26
$ws.Range("H26").Value = 7143.3335
$ws.Range("I26").Value = 7143.3335
$ws.Range("K26").Value = 7143.3335
$ws.Range("M26").Value = -6851.3335

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
# Row 15
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").ClearContents()
# Row 31
$ws.Range("H31").Value = 5937
$ws.Range("I31").Value = 7249.75
$ws.Range("J31").Value = 4624.25
$ws.Range("K31").Value = 7249.75
$ws.Range("L31").Value = 4624.25
$ws.Range("M31").Value = -6954.75
$ws.Range("N31").Value = -5214.25
# Row 34
$ws.Range("H34").Value = 5937
$ws.Range("I34").Value = 7249.75
$ws.Range("J34").Value = 4624.25
$ws.Range("K34").Value = 7249.75
$ws.Range("L34").Value = 4624.25
$ws.Range("M34").Value = -7047.75
$ws.Range("N34").Value = -5028.25
# Row 62
$ws.Range("H62").Value = 7380.2856
$ws.Range("I62").Value = 8763
$ws.Range("J62").Value = 6343.25
$ws.Range("K62").Value = 8763
$ws.Range("L62").Value = 6343.25
$ws.Range("M62").Value = -8139
$ws.Range("N62").Value = -7591.25
# Row 65
$ws.Range("H65").Value = 7380.2856
$ws.Range("I65").Value = 8763
$ws.Range("J65").Value = 6343.25
$ws.Range("K65").Value = 43815
$ws.Range("L65").Value = 31716.25
$ws.Range("M65").Value = -40695
$ws.Range("N65").Value = -37956.25

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
# Row 13
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("N13").Value = 0
$ws.Range("L13").ClearContents()
# Row 20
$ws.Range("H20").Value = 20
$ws.Range("I20").Value = 20
$ws.Range("K20").Value = 60
$ws.Range("M20").Value = 167
# Row 75
$ws.Range("H75").Value = 859.8
$ws.Range("I75").Value = 850
$ws.Range("J75").Value = 862.25
$ws.Range("K75").Value = 2550
$ws.Range("L75").Value = 2586.75
$ws.Range("M75").Value = -1552
$ws.Range("N75").Value = -4582.75
# Row 78
$ws.Range("H78").Value = 859.8
$ws.Range("I78").Value = 850
$ws.Range("J78").Value = 862.25
$ws.Range("K78").Value = 7650
$ws.Range("L78").Value = 7760.25
$ws.Range("M78").Value = -2658
$ws.Range("N78").Value = -17744.25
# Row 128
$ws.Range("H128").Value = 110000
$ws.Range("I128").Value = 110000
$ws.Range("K128").Value = 330000
$ws.Range("M128").Value = -325020
# Row 140
$ws.Range("H140").Value = 1961.8
$ws.Range("I140").Value = 1961.8
$ws.Range("K140").Value = 5885.4
$ws.Range("M140").Value = -705.3999999999996

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 4232.9
$ws.Range("I122").Value = 4257.8887
$ws.Range("K122").Value = 12773.6661
$ws.Range("M122").Value = -10323.6661
# Row 126
$ws.Range("H126").Value = 6500
$ws.Range("I126").Value = 7000
$ws.Range("J126").Value = 6000
$ws.Range("K126").Value = 21000
$ws.Range("L126").Value = 18000
$ws.Range("M126").Value = -18530
$ws.Range("N126").Value = -22940

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
# Row 17
$ws.Range("H17").Value = 5000
$ws.Range("I17").Value = 5000
$ws.Range("K17").Value = 5000
$ws.Range("M17").Value = -4830
# Row 22
$ws.Range("H22").Value = 1062
$ws.Range("I22").Value = 997.1667
$ws.Range("J22").Value = 1191.6666
$ws.Range("K22").Value = 997.1667
$ws.Range("L22").Value = 1191.6666
$ws.Range("M22").Value = -702.1667
$ws.Range("N22").Value = -1781.6666
# Row 27
$ws.Range("H27").Value = 1062
$ws.Range("I27").Value = 997.1667
$ws.Range("J27").Value = 1191.6666
$ws.Range("K27").Value = 997.1667
$ws.Range("L27").Value = 1191.6666
$ws.Range("M27").Value = -890.1667
$ws.Range("N27").Value = -1405.6666
# Row 40
$ws.Range("H40").Value = 5422.2144
$ws.Range("I40").Value = 4901.091
$ws.Range("K40").Value = 4901.091
$ws.Range("M40").Value = -4765.091
# Row 68
$ws.Range("H68").Value = 2353.818
$ws.Range("I68").Value = 2353.818
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 2353.818
$ws.Range("L68").Value = 0
$ws.Range("N68").Value = -1604.818
$ws.Range("M68").ClearContents()
# Row 71
$ws.Range("H71").Value = 2353.818
$ws.Range("I71").Value = 2353.818
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 11769.09
$ws.Range("L71").Value = 0
$ws.Range("N71").Value = -8025.09
$ws.Range("M71").ClearContents()
# Row 82
$ws.Range("H82").Value = 1876.5714
$ws.Range("I82").Value = 1843.6364
$ws.Range("J82").Value = 1997.3334
$ws.Range("K82").Value = 1843.6364
$ws.Range("L82").Value = 1997.3334
$ws.Range("M82").Value = -1482.6364
$ws.Range("N82").Value = -2719.3334
# Row 85
$ws.Range("H85").Value = 1876.5714
$ws.Range("I85").Value = 1843.6364
$ws.Range("J85").Value = 1997.3334
$ws.Range("K85").Value = 1843.6364
$ws.Range("L85").Value = 1997.3334
$ws.Range("M85").Value = -595.6364000000001
$ws.Range("N85").Value = -4493.3334
# Row 122
$ws.Range("H122").Value = 3606.4546
$ws.Range("I122").Value = 3606.4546
$ws.Range("K122").Value = 10819.3638
$ws.Range("M122").Value = -8369.363799999999
# Row 132
$ws.Range("H132").Value = 4918.7
$ws.Range("I132").Value = 4918.7
$ws.Range("K132").Value = 14756.1
$ws.Range("M132").Value = -12226.1
# Row 136
$ws.Range("H136").Value = 2500
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
# Row 12
$ws.Range("H12").Value = 8333
$ws.Range("I12").Value = 7999.5
$ws.Range("K12").Value = 7999.5
$ws.Range("M12").Value = -7857.5
